{"js": "// Reorders the \"COMPETENCES TECHNIQUES\" skill lines:\n//   before: Langages, Bases de donn\u00e9es, Autres, Visualisation, ML/AI, MLOps\n//   after : Langages, Visualisation, MLOps, Autres, ML/AI, Bases de donn\u00e9es\n// (the resume_data source was converted from .py to .json, changing dict/key\n// iteration order \u2014 the wording of every line is unchanged, only their\n// paragraph order is).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The block of 6 skill lines immediately follows the \"Langages : ...\" line.\nconst targetOrder = [\n  \"Langages : r, python, matlab, c, c++\",\n  \"Visualisation : web analytics, tableau\",\n  \"MLOps : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\",\n  \"Autres : marketing, google analytics, internes comme externes, presse, affichage, site centric, formats\",\n  \"ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n  \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\",\n];\n\nlet startIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Langages : r, python\") === 0) {\n    startIndex = i;\n    break;\n  }\n}\nif (startIndex === -1) {\n  throw new Error(\"Could not locate the 'Langages : ...' skills paragraph.\");\n}\n\n// Paragraph 0 of the block (\"Langages\") is unchanged; only replace the text\n// of the next 5 paragraphs with the reordered content so every paragraph\n// keeps its original formatting (pPr/rPr) and only the run text moves.\nfor (let k = 1; k < targetOrder.length; k++) {\n  const paragraph = paragraphs.items[startIndex + k];\n  const range = paragraph.getRange();\n  range.insertText(targetOrder[k], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Reorders the \"COMPETENCES TECHNIQUES\" skill lines:\n#   before: Langages, Bases de donnees, Autres, Visualisation, ML/AI, MLOps\n#   after : Langages, Visualisation, MLOps, Autres, ML/AI, Bases de donnees\n# (the resume_data source was converted from .py to .json, changing dict/key\n# iteration order -- the wording of every line is unchanged, only their\n# paragraph order is).\n\n$d = $word.ActiveDocument\n\n$targetOrder = @(\n    \"Langages : r, python, matlab, c, c++\",\n    \"Visualisation : web analytics, tableau\",\n    \"MLOps : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\",\n    \"Autres : marketing, google analytics, internes comme externes, presse, affichage, site centric, formats\",\n    \"ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n    \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\"\n)\n\n$count = $d.Paragraphs.Count\n$startIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.StartsWith(\"Langages : r, python\")) {\n        $startIdx = $i\n        break\n    }\n}\nif ($startIdx -eq -1) {\n    throw \"Could not locate the 'Langages : ...' skills paragraph.\"\n}\n\n# Paragraph 0 of the block (\"Langages\") is unchanged; only replace the text\n# of the next 5 paragraphs with the reordered content so every paragraph\n# keeps its original formatting and only the run text moves.\nfor ($k = 1; $k -lt $targetOrder.Length; $k++) {\n    $p = $d.Paragraphs.Item($startIdx + $k)\n    $p.Range.Text = $targetOrder[$k]\n}\n"}
